# Applies the commit's change set to PlayerPerformance_5664.xlsx:
#   1. Removes the stray empty inline-string cell at "ODI Batting"!B2.
#   2. Adds a new trailing worksheet "ODI Batting Extra" (sheetId 4) with a
#      bold header row and two data rows, mirroring the header style already
#      used on the other sheets.

$wb = $excel.ActiveWorkbook

# Helper: write a value that must be stored as TEXT (not auto-coerced to a
# number/percentage/date by Excel's normal input parsing), then drop the
# temporary "@" text format again so the cell doesn't pick up a stray style
# index that wasn't present in the source data.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- 1. "ODI Batting" : drop the empty B2 cell -----------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("B2").ClearContents()

# --- 2. Add "ODI Batting Extra" as the new last sheet ----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExtra = $wb.Worksheets.Add($null, $lastSheet)
$wsExtra.Name = "ODI Batting Extra"

# Reuse the same bold/bordered/centered header style already used by the
# other sheets (e.g. "ODI Batting"!A1:F1) instead of building a new one.
$wsBatting.Range("A1:F1").Copy()
$wsExtra.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"

# Row 2
Set-TextValue $wsExtra.Range("A2") "4433"
Set-TextValue $wsExtra.Range("B2") ""
Set-TextValue $wsExtra.Range("C2") ""
Set-TextValue $wsExtra.Range("D2") ""
Set-TextValue $wsExtra.Range("E2") ""
Set-TextValue $wsExtra.Range("F2") "NO"

# Row 3
Set-TextValue $wsExtra.Range("A3") "4434"
$wsExtra.Range("B3").Value = 10
Set-TextValue $wsExtra.Range("C3") "2"
Set-TextValue $wsExtra.Range("D3") "0"
Set-TextValue $wsExtra.Range("E3") "3.24%"
Set-TextValue $wsExtra.Range("F3") "NO"
